# Populate column H (Срок изготовления партии, дней) for rows 2-6
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 457
$ws.Range("H3").Value = 456
$ws.Range("H4").Value = 4380
$ws.Range("H5").Value = 234
$ws.Range("H6").Value = 9875
